$d = $word.ActiveDocument

# -----------------------------------------------------------------
# Edit 1: the empty table cell in the "week" schedule table (row for
# "1 / Fri, Nov 1, 13 / Intro to visuality / <empty>") should get the
# same "Compact" paragraph style used by its sibling cells.
# -----------------------------------------------------------------
$table = $d.Tables.Item(1)
for ($r = 1; $r -le $table.Rows.Count; $r++) {
    for ($c = 1; $c -le $table.Columns.Count; $c++) {
        $cell = $null
        try {
            $cell = $table.Cell($r, $c)
        } catch {
            $cell = $null
        }
        if ($cell -ne $null) {
            $para = $cell.Range.Paragraphs.Item(1)
            # A cell's/paragraph's Range.Text always carries a trailing
            # cell-/paragraph-mark (CR, and CR+BEL for the last cell in
            # a row) - strip those control characters before testing
            # whether the paragraph has any real text in it.
            $visibleText = $para.Range.Text.Trim([char]13, [char]7)
            if ($visibleText -eq "") {
                $para.Style = "Compact"
            }
        }
    }
}

# -----------------------------------------------------------------
# Edit 2: rename the auto-generated hash bookmark around the
# "Visuality in Educational Media or Youth Media (50%)" Heading3 to a
# human-readable slug. This runtime's Bookmark.Name setter / Delete()
# are both no-ops, so rebuild the bookmark: delete the old
# heading paragraph (a zero-width delete positioned at the bookmark's
# start removes the bookmarked paragraph + its bookmark together),
# re-insert an equivalent Heading3 paragraph with the same text, and
# add a fresh bookmark with the desired name around that text.
# -----------------------------------------------------------------
$headingText = "Visuality in Educational Media or Youth Media (50%)"
$oldBookmarkName = "Xdd32528e7634130f5d24f03120080f817a9565c"
$newBookmarkName = "visuality-in-educational-media-or-youth-media-50"

$targetBookmark = $null
for ($i = 1; $i -le $d.Bookmarks.Count; $i++) {
    $b = $d.Bookmarks.Item($i)
    if ($b.Name -eq $oldBookmarkName) {
        $targetBookmark = $b
        break
    }
}

if ($targetBookmark -ne $null) {
    $pos = $targetBookmark.Start

    # Remove the old heading paragraph (text + its bookmark) in one go.
    $zeroRange = $d.Range($pos, $pos)
    $zeroRange.Delete()

    # Re-insert a fresh paragraph at the same position.
    $insertRange = $d.Range($pos, $pos)
    $insertRange.InsertParagraphAfter()

    # Find that new (currently empty) paragraph and give it the
    # Heading3 style, matching the original.
    $newPara = $null
    for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
        $p = $d.Paragraphs.Item($i)
        if ($p.Range.Start -eq $pos) {
            $newPara = $p
            break
        }
    }
    $newPara.Style = "Heading3"

    # Insert the heading text.
    $textInsertPoint = $d.Range($pos, $pos)
    $textInsertPoint.InsertAfter($headingText)

    # Wrap the new text in a bookmark with the new slug name.
    $textEnd = $pos + $headingText.Length
    $bmRange = $d.Range($pos, $textEnd)
    $d.Bookmarks.Add($newBookmarkName, $bmRange)
}
